# Update gh-pages to output generated at 456a3b4
# Apply updated "想去人数" (F) and "最低票价" (G) values to the
# "展览" and "全部类型" worksheets (which hold identical data tables).

$wb = $excel.ActiveWorkbook

# Row number -> hashtable of column letter -> new value
$updates = @{
    2  = @{ F = 699 }
    4  = @{ F = 247 }
    9  = @{ F = 5974 }
    11 = @{ F = 326 }
    12 = @{ F = 262; G = 55 }
    16 = @{ F = 4709 }
    18 = @{ F = 1247 }
    20 = @{ F = 99 }
    21 = @{ F = 215 }
    25 = @{ F = 157 }
    27 = @{ F = 370 }
    28 = @{ F = 55 }
    29 = @{ F = 49 }
    31 = @{ F = 31 }
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $cols = $updates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
